$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New IPC PO (column C) values for rows 2-51, from updated model predictions
$newC = @{}
$newC[2] = 29.58759074234123
$newC[3] = 30.18283030166088
$newC[4] = 30.05240729647354
$newC[5] = 29.8008296316535
$newC[6] = 29.81086508372359
$newC[7] = 29.8588191652655
$newC[8] = 29.96757500962112
$newC[9] = 30.10526514425146
$newC[10] = 30.16667345311208
$newC[11] = 30.45523101141733
$newC[12] = 30.47076606092922
$newC[13] = 30.43801446389514
$newC[14] = 30.69975840814894
$newC[15] = 30.82994133128413
$newC[16] = 30.96543978898074
$newC[17] = 31.14081628343929
$newC[18] = 31.21232837632543
$newC[19] = 31.04908220229241
$newC[20] = 30.83156346787992
$newC[21] = 31.05406723183888
$newC[22] = 31.59665854667251
$newC[23] = 32.49295279257779
$newC[24] = 32.51761548091525
$newC[25] = 32.75503165065629
$newC[26] = 32.96116689884062
$newC[27] = 33.02125648673261
$newC[28] = 33.17862983298771
$newC[29] = 33.54432577702764
$newC[30] = 33.59484823058081
$newC[31] = 33.6974953360889
$newC[32] = 34.2781710818468
$newC[33] = 34.59026156137988
$newC[34] = 35.69277273257615
$newC[35] = 35.95596127809328
$newC[36] = 36.34356418676805
$newC[37] = 36.87894717070319
$newC[38] = 37.02797098788045
$newC[39] = 37.90215637364533
$newC[40] = 38.58213526225354
$newC[41] = 39.10436487970355
$newC[42] = 39.423575650698
$newC[43] = 39.69720863308636
$newC[44] = 39.75404246321752
$newC[45] = 39.89197518924603
$newC[46] = 40.18348844285919
$newC[47] = 41.19773897167246
$newC[48] = 41.27961779414809
$newC[49] = 41.72409164356247
$newC[50] = 42.39997101360544
$newC[51] = 43.53406211122057

$deltaSum = 0
$deltaSqSum = 0

for ($r = 2; $r -le 51; $r++) {
    $b = $ws.Cells.Item($r, 2).Value2
    $c = $newC[$r]
    $ws.Cells.Item($r, 3).Value2 = $c
    $d = $c - $b
    $e = $d * $d
    $ws.Cells.Item($r, 4).Value2 = $d
    $ws.Cells.Item($r, 5).Value2 = $e
    $deltaSum = $deltaSum + $d
    $deltaSqSum = $deltaSqSum + $e
}

# Row 52: TOTAL (sum of DELTA in C52, sum of DELTA^2 in E52)
$ws.Cells.Item(52, 3).Value2 = $deltaSum
$ws.Cells.Item(52, 5).Value2 = $deltaSqSum

# Row 53: MSE (average of DELTA^2)
$mse = $deltaSqSum / 50
$ws.Cells.Item(53, 5).Value2 = $mse

Write-Host "Applied sliding window results update"
